$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.890.12'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '1.638.05'
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '''215.29'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").Value = '''0.5095'
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = '''1.004'
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").Value = '''0.2587'
$ws.Range("E8").Value = '  +0.72%  '
$ws.Range("D9").Value = '''0.06435'
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("D10").Value = '''20.35'
$ws.Range("E10").Value = '  +4.68%  '
$ws.Range("D11").Value = '''0.07796'
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.656.96'
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.281'
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.863.76'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '''0.5603'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").Value = '0.0₅7678'
$ws.Range("E16").Value = '  +2.51%  '
$ws.Range("D17").Value = '''63.25'
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("D18").Value = '25.886.84'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").Value = '''193.99'
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("D21").Value = '''4.389'
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("D22").Value = '''9.965'
$ws.Range("E22").Value = '  +1.91%  '
$ws.Range("D23").Value = '''6.158'
$ws.Range("E23").Value = '  +2.37%  '
$ws.Range("D24").Value = '''1.004'
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("D25").Value = '''1.786'
$ws.Range("E25").Value = '  -4.95%  '
$ws.Range("D26").Value = '''138.26'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").Value = '''6.842'
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("D29").Value = '''15.61'
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = '''0.04972'
$ws.Range("E31").Value = '  +2.12%  '
$ws.Range("D32").Value = '''3.302'
$ws.Range("E32").Value = '  +1.58%  '
$ws.Range("D33").Value = '''3.256'
$ws.Range("E33").Value = '  +2.61%  '
$ws.Range("D34").Value = '''1.571'
$ws.Range("E34").Value = '  +1.88%  '
$ws.Range("D35").Value = '''2.388'
$ws.Range("E35").Value = '  +1.18%  '
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("D37").Value = '''2.579'
$ws.Range("E37").Value = '  +1.55%  '
$ws.Range("D38").Value = '''0.5572'
$ws.Range("D39").Value = '1.136.86'
$ws.Range("E39").Value = '  +1.98%  '
$ws.Range("D40").Value = '''0.01575'
$ws.Range("E40").Value = '  +1.75%  '
$ws.Range("D41").Value = '''0.9975'
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E42").Value = '  +2.42%  '
$ws.Range("D43").Value = '''5.476'
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").Value = '''0.8032'
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("D45").Value = '0.0₈111'
$ws.Range("E45").Value = '  -4.03%  '
$ws.Range("D46").Value = '''55.57'
$ws.Range("E46").Value = '  +1.72%  '
$ws.Range("D47").Value = '''0.4266'
$ws.Range("E47").Value = '  -3.55%  '
$ws.Range("D48").Value = '''7.784'
$ws.Range("E48").Value = '  +2.82%  '
$ws.Range("D49").Value = '''0.05073'
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("D51").Value = '''1.003'
$ws.Range("E51").Value = '  +0.36%  '
